$wb = $excel.ActiveWorkbook

# Update the "Ready for handoff" status everywhere it is shown: the summary Overview sheet
# (columns B/C) as well as the per-locale Status column (B) on the zh-cn/de-de sheets.
$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("B3").Value = "Handed back: in sync with en-US"

# zh-cn: populate "Latest Target File" (E) / "Latest Handback File" (F) for the two handed-back
# rows, and refresh "Latest Handback DateTime" (G).
$zh.Range("E2").Value = "0eaf48a6-a0f0-4054-bb50-ed938ed0e546.md"
$zh.Range("F2").Value = "0eaf48a6-a0f0-4054-bb50-ed938ed0e546.b68dcd8ba7a6f8da81393ba802d82e6837f4ee46.zh-cn.xlf"
$zh.Range("G2").Value = "2016-03-08 03:03:16"

$zh.Range("E3").Value = "b33877fe-e3c8-4168-bf9f-da162a11ce8a.md"
$zh.Range("F3").Value = "b33877fe-e3c8-4168-bf9f-da162a11ce8a.c478d0514c0d914a3a5319b97bdfe4331471b48f.zh-cn.xlf"
$zh.Range("G3").Value = "2016-03-08 03:03:16"

$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/835bcf548ec8b7e04d92253ce8e6bd7a4480e38f/e2e/0eaf48a6-a0f0-4054-bb50-ed938ed0e546.md", [System.Type]::Missing, [System.Type]::Missing, "0eaf48a6-a0f0-4054-bb50-ed938ed0e546.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09bc4358632257fffbc22148e7a53dedb66ada33/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/high/0eaf48a6-a0f0-4054-bb50-ed938ed0e546.b68dcd8ba7a6f8da81393ba802d82e6837f4ee46.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "0eaf48a6-a0f0-4054-bb50-ed938ed0e546.b68dcd8ba7a6f8da81393ba802d82e6837f4ee46.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/835bcf548ec8b7e04d92253ce8e6bd7a4480e38f/e2e/b33877fe-e3c8-4168-bf9f-da162a11ce8a.md", [System.Type]::Missing, [System.Type]::Missing, "b33877fe-e3c8-4168-bf9f-da162a11ce8a.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09bc4358632257fffbc22148e7a53dedb66ada33/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/high/b33877fe-e3c8-4168-bf9f-da162a11ce8a.c478d0514c0d914a3a5319b97bdfe4331471b48f.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "b33877fe-e3c8-4168-bf9f-da162a11ce8a.c478d0514c0d914a3a5319b97bdfe4331471b48f.zh-cn.xlf")

# de-de: same pattern, with the de-de handback file names / timestamp.
$de.Range("E2").Value = "0eaf48a6-a0f0-4054-bb50-ed938ed0e546.md"
$de.Range("F2").Value = "0eaf48a6-a0f0-4054-bb50-ed938ed0e546.b68dcd8ba7a6f8da81393ba802d82e6837f4ee46.de-de.xlf"
$de.Range("G2").Value = "2016-03-08 03:03:33"

$de.Range("E3").Value = "b33877fe-e3c8-4168-bf9f-da162a11ce8a.md"
$de.Range("F3").Value = "b33877fe-e3c8-4168-bf9f-da162a11ce8a.c478d0514c0d914a3a5319b97bdfe4331471b48f.de-de.xlf"
$de.Range("G3").Value = "2016-03-08 03:03:33"

$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/835bcf548ec8b7e04d92253ce8e6bd7a4480e38f/e2e/0eaf48a6-a0f0-4054-bb50-ed938ed0e546.md", [System.Type]::Missing, [System.Type]::Missing, "0eaf48a6-a0f0-4054-bb50-ed938ed0e546.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff24cb31c23f30df06e5bf43eec74e027e5cd2e1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/high/0eaf48a6-a0f0-4054-bb50-ed938ed0e546.b68dcd8ba7a6f8da81393ba802d82e6837f4ee46.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "0eaf48a6-a0f0-4054-bb50-ed938ed0e546.b68dcd8ba7a6f8da81393ba802d82e6837f4ee46.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/835bcf548ec8b7e04d92253ce8e6bd7a4480e38f/e2e/b33877fe-e3c8-4168-bf9f-da162a11ce8a.md", [System.Type]::Missing, [System.Type]::Missing, "b33877fe-e3c8-4168-bf9f-da162a11ce8a.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff24cb31c23f30df06e5bf43eec74e027e5cd2e1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/high/b33877fe-e3c8-4168-bf9f-da162a11ce8a.c478d0514c0d914a3a5319b97bdfe4331471b48f.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "b33877fe-e3c8-4168-bf9f-da162a11ce8a.c478d0514c0d914a3a5319b97bdfe4331471b48f.de-de.xlf")
